# act tablas web jul25
$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# --- Data sheet: add two new years (2023, 2022) at the top of the series ---
$wsData = $wb.Worksheets.Item("Data")

# insert two new rows right after the header row, pushing existing data down
$wsData.Range("A2:A3").EntireRow.Insert()

Set-TextValue $wsData.Cells.Item(2,1) "2023"
$wsData.Cells.Item(2,2).Value = 4.1
Set-TextValue $wsData.Cells.Item(3,1) "2022"
$wsData.Cells.Item(3,2).Value = 4.1

# --- Metadata sheet: fix blank cell + add "actualizacion" row before "cita" ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# A1 was an empty string; it should now contain a single space
Set-TextValue $wsMeta.Cells.Item(1,1) " "

# insert a new row before the "cita" row (currently row 9)
$wsMeta.Rows.Item(9).Insert()
Set-TextValue $wsMeta.Cells.Item(9,1) "actualizacion"
Set-TextValue $wsMeta.Cells.Item(9,2) "Julio 2025"
